# Fruta / hortaliza, semanal
# Insert a new week of price observations (2 rows: "Primera" and "Segunda"
# quality for Packham's Triumph) at the top of the Packham's Triumph /
# Winter Nelis / Abate Fettel rotating block, pushing all the existing
# observations down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (rows 237-266) down by two rows to make room for the
# newest week's observations.
$ws.Rows("237:238").Insert()

# Row 237: Packham's Triumph - Primera
$ws.Cells.Item(237, 1).Value = 11
$ws.Cells.Item(237, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(237, 3).Value = "Bíobío"
$ws.Cells.Item(237, 4).Value = 44474
$ws.Cells.Item(237, 5).Value = 8
$ws.Cells.Item(237, 6).Value = "Fruta"
$ws.Cells.Item(237, 7).Value = 100104
$ws.Cells.Item(237, 8).Value = "Frutos de pepita"
$ws.Cells.Item(237, 9).Value = 100104005
$ws.Cells.Item(237, 10).Value = "Pera"
$ws.Cells.Item(237, 11).Value = "Packham's Triumph"
$ws.Cells.Item(237, 12).Value = "Primera"
$ws.Cells.Item(237, 13).Value = 100
$ws.Cells.Item(237, 14).Value = 9000
$ws.Cells.Item(237, 15).Value = 10000
$ws.Cells.Item(237, 16).Value = 9500
$ws.Cells.Item(237, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(237, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(237, 19).Value = 594
$ws.Cells.Item(237, 20).Value = 16

# Row 238: Packham's Triumph - Segunda
$ws.Cells.Item(238, 1).Value = 11
$ws.Cells.Item(238, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(238, 3).Value = "Bíobío"
$ws.Cells.Item(238, 4).Value = 44474
$ws.Cells.Item(238, 5).Value = 8
$ws.Cells.Item(238, 6).Value = "Fruta"
$ws.Cells.Item(238, 7).Value = 100104
$ws.Cells.Item(238, 8).Value = "Frutos de pepita"
$ws.Cells.Item(238, 9).Value = 100104005
$ws.Cells.Item(238, 10).Value = "Pera"
$ws.Cells.Item(238, 11).Value = "Packham's Triumph"
$ws.Cells.Item(238, 12).Value = "Segunda"
$ws.Cells.Item(238, 13).Value = 50
$ws.Cells.Item(238, 14).Value = 8000
$ws.Cells.Item(238, 15).Value = 8000
$ws.Cells.Item(238, 16).Value = 8000
$ws.Cells.Item(238, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(238, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(238, 19).Value = 500
$ws.Cells.Item(238, 20).Value = 16
